$d = $word.ActiveDocument

function Insert-Break {
    param(
        [string]$needle
    )
    $replacement = "^l" + $needle
    $d.Content.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $replacement, 2) | Out-Null
}

# --- "Programa resumido" paragraph ---
Insert-Break "2.Composição dos vidros"
Insert-Break "3.Materiais Primas"
Insert-Break "4.Mecanismo de fusão e formação do vidro"
Insert-Break "5.Viscosidade - Definição"
Insert-Break "6.Propriedades óticas "
Insert-Break "7.Propriedades mecânicas "
Insert-Break "8.Propriedades químicas "
Insert-Break "9.Processamento - Vidro plano"
Insert-Break "10.Aula prática - Fundir um vidro"

# --- "Programa" paragraph ---
Insert-Break "2. Composição dos vidros"
Insert-Break "3. Materiais Primas"
Insert-Break "4. Mecanismo de fusão e formação do vidro"
Insert-Break "5. Viscosidade – Definição"
Insert-Break "6. Propriedades óticas "
Insert-Break "7. Propriedades mecânicas "
Insert-Break "8. Propriedades químicas "
Insert-Break "9. Processamento – Vidro plano"
Insert-Break "10. Aula prática - Fundir um vidro"

# --- "Bibliografia" paragraph ---
Insert-Break "2.)H. Scholze, Glas, Springer-Verlag, 1988"
Insert-Break "3.)R. H. Doremus, Glass Science, New York, John Wiley, 1994"
Insert-Break "4.)H. G. Pfaender, Schott Guide to Glass, London, Chapman & Hall, 1996"
